$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (fix typo introduced by author)
$ws.Name = 'Transaltions'

# Populate translation rows (reordered + new entries)
$ws.Cells.Item(1,1).Value = 'Name'
$ws.Cells.Item(1,2).Value = 'Source Language'
$ws.Cells.Item(1,3).Value = 'Target Language'
$ws.Cells.Item(1,4).Value = 'Status'

$ws.Cells.Item(2,1).Value = 'Hello'
$ws.Cells.Item(2,2).Value = 'Hello'
$ws.Cells.Item(2,3).Value = 'Hallå'
$ws.Cells.Item(2,4).Value = 'Need review'

$ws.Cells.Item(3,1).Value = 'Press_Any_Key_To_Exit'
$ws.Cells.Item(3,2).Value = 'Press any key to exit'
$ws.Cells.Item(3,3).Value = 'tryck på valfri tangent för att avsluta'
$ws.Cells.Item(3,4).Value = 'Need review'

$ws.Cells.Item(4,1).Value = 'Enter_Email'
$ws.Cells.Item(4,2).Value = 'Enter email please!'
$ws.Cells.Item(4,3).Value = 'Ange email tack!'
$ws.Cells.Item(4,4).Value = 'Need review'

$ws.Cells.Item(5,1).Value = 'Enter_Email_To_Be_The_Next_User_To_Be_Able_To_Login_To_The_Best_Office_System_In_The_World'
$ws.Cells.Item(5,2).Value = 'Enter email please!'
$ws.Cells.Item(5,3).Value = 'Ange email tack!'
$ws.Cells.Item(5,4).Value = 'Need review'

$ws.Cells.Item(6,1).Value = 'Password'
$ws.Cells.Item(6,2).Value = 'Password'
$ws.Cells.Item(6,3).Value = 'Lösenord'
$ws.Cells.Item(6,4).Value = 'Need review'

$ws.Cells.Item(7,1).Value = 'Home'
$ws.Cells.Item(7,2).Value = 'Home'
$ws.Cells.Item(7,3).Value = 'Hem'
$ws.Cells.Item(7,4).Value = 'Need review'

$ws.Cells.Item(8,1).Value = 'Dog'
$ws.Cells.Item(8,2).Value = 'Dog'
$ws.Cells.Item(8,3).Value = 'Hund'
$ws.Cells.Item(8,4).Value = 'Final'

$ws.Cells.Item(9,1).Value = 'Cat'
$ws.Cells.Item(9,2).Value = 'Cat'
$ws.Cells.Item(9,3).Value = 'Katt'
$ws.Cells.Item(9,4).Value = 'Need review'

$ws.Cells.Item(10,1).Value = 'Bird'
$ws.Cells.Item(10,2).Value = 'Bird'
$ws.Cells.Item(10,3).Value = 'Fågel'
$ws.Cells.Item(10,4).Value = 'Need review'

$ws.Cells.Item(11,1).Value = 'Cow'
$ws.Cells.Item(11,2).Value = 'Cow'
$ws.Cells.Item(11,3).Value = 'Ko'
$ws.Cells.Item(11,4).Value = 'Need review'

$ws.Cells.Item(12,1).Value = 'Love_Is'
$ws.Cells.Item(12,2).Value = 'Love is'
$ws.Cells.Item(12,3).Value = 'Kärlek är'
$ws.Cells.Item(12,4).Value = 'Need review'

$ws.Cells.Item(13,1).Value = 'House'
$ws.Cells.Item(13,2).Value = 'House'
$ws.Cells.Item(13,3).Value = 'Hus'
$ws.Cells.Item(13,4).Value = 'Need review'

$ws.Cells.Item(14,1).Value = 'Street'
$ws.Cells.Item(14,2).Value = 'Street'
$ws.Cells.Item(14,3).Value = 'Gata'
$ws.Cells.Item(14,4).Value = 'Final'

$ws.Cells.Item(15,1).Value = 'Street_And_House2'
$ws.Cells.Item(15,2).Value = 'Street and House and Room. later this day i was so sad that i jumped'
$ws.Cells.Item(15,3).Value = 'Gata och Hus och Rum. senare denna dag var jag så ledsen att jag hoppade'
$ws.Cells.Item(15,4).Value = 'Need review'

$ws.Cells.Item(16,1).Value = 'Computer'
$ws.Cells.Item(16,2).Value = 'Computer'
$ws.Cells.Item(16,3).Value = 'Dator'
$ws.Cells.Item(16,4).Value = 'Need review'

$ws.Cells.Item(17,1).Value = 'Cup'
$ws.Cells.Item(17,2).Value = 'Cup'
$ws.Cells.Item(17,3).Value = 'Kopp'
$ws.Cells.Item(17,4).Value = 'Need review'

$ws.Cells.Item(18,1).Value = 'Read_Instructions'
$ws.Cells.Item(18,2).Value = 'Read
Instructions'
$ws.Cells.Item(18,3).Value = 'Läsa
Instruktioner'
$ws.Cells.Item(18,4).Value = 'Need review'

$ws.Cells.Item(19,1).Value = 'Nice_Gool_Nice_Shot'
$ws.Cells.Item(19,2).Value = 'Nice goal. Nice Shot.'
$ws.Cells.Item(19,3).Value = 'Snyggt mål. Snyggt skott.'
$ws.Cells.Item(19,4).Value = 'Need review'

# Bold header row
$ws.Range("A1:D1").Font.Bold = $true

# Column widths (best effort to match source widths)
$ws.Columns.Item(1).ColumnWidth = 89.83333333333333
$ws.Columns.Item(2).ColumnWidth = 68.83333333333333
$ws.Columns.Item(3).ColumnWidth = 71.83333333333333
$ws.Columns.Item(4).ColumnWidth = 10.833333333333334

# Row heights
for ($r = 1; $r -le 19; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}

# Restore selected cell
[void]$ws.Range("D13").Select()
